$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '322.17'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '8.08%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '48.56'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '15.08%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.279'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '5.40%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08115'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '7.91%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.582'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '4.93%'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '3.19%'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '29.62%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1299'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '9.83%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1948'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '5.95%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09502'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '6.26%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04630'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '12.06%'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.02%'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '3.72%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005848'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.38%'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.14%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.426'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.02%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '2.15%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.091'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-3.01%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1410'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '4.34%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.3125'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.67%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04257'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '4.24%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001305'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '3.05%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004251'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '9.31%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003540'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-4.94%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02695'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '12.38%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05667'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '8.48%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.006301'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-7.21%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007676'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.15%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1439'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '8.72%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.007696'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '3.78%'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '13.89%'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '6.77%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00007006'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '6.56%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.07%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05405'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-1.50%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.004001'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-4.83%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.07%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.07%'
